# Trade #8 closed at 2026-02-16 21:51:45 - momentum DOWN +0.000%
#
# The bot's "momentum" strategy fired for the first time, so a brand-new
# "momentum" tab is introduced (taking over the old "leadlag" sheetId/rId),
# a fresh copy of "leadlag" is appended right after it to keep the existing
# leadlag trade log intact, and the new trade row is appended to both the
# new "momentum" sheet and the master "All Trades" log.

$wb = $excel.ActiveWorkbook

# --- Step 1: duplicate "leadlag" so its trade history survives unchanged.
# Placing the copy immediately after the original gives us, in tab order:
#   Summary, Strategy Status, All Trades, leadlag, leadlag (2)
$leadlagSheet = $wb.Worksheets.Item("leadlag")
$leadlagSheet.Copy([System.Reflection.Missing]::Value, $leadlagSheet)
$newCopy = $wb.Worksheets.Item("leadlag (2)")

# --- Step 2: the ORIGINAL "leadlag" sheet is repurposed as "momentum"
# (keeps its original sheetId/r:id); the copy becomes the "new" leadlag
# sheet (gets a freshly minted sheetId/r:id). Final tab order:
#   Summary, Strategy Status, All Trades, momentum, leadlag
$leadlagSheet.Name = "momentum"
$newCopy.Name = "leadlag"

# --- Step 3: the momentum sheet should contain only the header row plus
# this one new trade, so drop the 7 leadlag rows it inherited from the copy.
$momentum = $wb.Worksheets.Item("momentum")
$momentum.Rows("3:8").Delete()

$momentum.Range("A2").Value = 8
$momentum.Range("B2").Value = "'2026-02-16"
$momentum.Range("C2").Value = "21:51:45"
$momentum.Range("D2").Value = "momentum"
$momentum.Range("E2").Value = "DOWN"
$momentum.Range("F2").Value = 68234.395
$momentum.Range("G2").Value = "'"
$momentum.Range("H2").Value = "OPEN"
$momentum.Range("I2").Value = 0
$momentum.Range("J2").Value = 0
$momentum.Range("K2").Value = 100
$momentum.Range("L2").Value = 0.9
$momentum.Range("M2").Value = "Downward momentum: -0.197% over 10 samples"
$momentum.Range("N2").Value = "'"
$momentum.Range("O2").Value = 0

# The leading "'" on B2/G2/N2 forces text (so the date string isn't parsed
# into a serial number, and the blank exit/exit-reason stay empty TEXT
# cells rather than being cleared to blank) without leaving Excel's
# "quote prefix" number format attached to the cell afterwards.
$momentum.Range("B2").Style = "Normal"
$momentum.Range("G2").Style = "Normal"
$momentum.Range("N2").Style = "Normal"

# --- Step 4: mirror the same trade as row 9 of the master "All Trades" log.
$allTrades = $wb.Worksheets.Item("All Trades")
$allTrades.Range("A9").Value = 8
$allTrades.Range("B9").Value = "'2026-02-16"
$allTrades.Range("C9").Value = "21:51:45"
$allTrades.Range("D9").Value = "momentum"
$allTrades.Range("E9").Value = "DOWN"
$allTrades.Range("F9").Value = 68234.395
$allTrades.Range("G9").Value = "'"
$allTrades.Range("H9").Value = "OPEN"
$allTrades.Range("I9").Value = 0
$allTrades.Range("J9").Value = 0
$allTrades.Range("K9").Value = 100
$allTrades.Range("L9").Value = 0.9
$allTrades.Range("M9").Value = "Downward momentum: -0.197% over 10 samples"
$allTrades.Range("N9").Value = "'"
$allTrades.Range("O9").Value = 0

$allTrades.Range("B9").Style = "Normal"
$allTrades.Range("G9").Style = "Normal"
$allTrades.Range("N9").Style = "Normal"
